$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 18 - Belgica
$ws.Range("B18").Value = 53449
$ws.Range("C18").Value = 368
$ws.Range("D18").Value = 13697
$ws.Range("E18").Value = 31045
$ws.Range("F18").Value = 478
$ws.Range("G18").Value = 51
$ws.Range("H18").Value = 8707

# Row 39 - Rumania
$ws.Range("E39").Value = 7339
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 972

# Row 57 - Finlandia
$ws.Range("B57").Value = 5984
$ws.Range("C57").Value = 22
$ws.Range("E57").Value = 1717

# Row 63 - Afganistan
$ws.Range("B63").Value = 4687
$ws.Range("C63").Value = 285
$ws.Range("E63").Value = 4007
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 122

# Row 75 - Uzbekistan
$ws.Range("D75").Value = 1945
$ws.Range("E75").Value = 498

# Row 129 - Estado de Palestina
$ws.Range("D129").Value = 301
$ws.Range("E129").Value = 72
